$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.090.09"
$ws.Range("E2").Value = "  +11.39%  "
$ws.Range("D3").Value = "'1.815.05"
$ws.Range("E3").Value = "  +7.97%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'228.49"
$ws.Range("E5").Value = "  +3.61%  "
$ws.Range("D6").Value = "'0.545"
$ws.Range("E6").Value = "  +2.98%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'31.38"
$ws.Range("E8").Value = "  +3.74%  "
$ws.Range("D9").Value = "'46.69"
$ws.Range("E9").Value = "  +5.31%  "
$ws.Range("E10").Value = "  +6.06%  "
$ws.Range("D11").Value = "'0.0667"
$ws.Range("E11").Value = "  +2.67%  "
$ws.Range("D12").Value = "'0.0928"
$ws.Range("E12").Value = "  +2.33%  "
$ws.Range("D13").Value = "'2.073.95"
$ws.Range("E13").Value = "  +8.00%  "
$ws.Range("D14").Value = "'1.812.03"
$ws.Range("E14").Value = "  +7.96%  "
$ws.Range("D15").Value = "'0.642"
$ws.Range("E15").Value = "  +3.67%  "
$ws.Range("D16").Value = "'34.052.60"
$ws.Range("E16").Value = "  +11.19%  "
$ws.Range("D17").Value = "'10.27"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").Value = "'4.26"
$ws.Range("E18").Value = "  +6.86%  "
$ws.Range("D19").Value = "'69.44"
$ws.Range("E19").Value = "  +4.40%  "
$ws.Range("D20").Value = "'258.01"
$ws.Range("E20").Value = "  +4.96%  "
$ws.Range("D21").Value = "'0.0₃0747"
$ws.Range("E21").Value = "  +2.57%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").Value = "'10.50"
$ws.Range("E23").Value = "  +3.89%  "
$ws.Range("D24").Value = "'4.34"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("E25").Value = "  +2.10%  "
$ws.Range("D26").Value = "'157.94"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "'16.57"
$ws.Range("E27").Value = "  +3.82%  "
$ws.Range("D28").Value = "'7.13"
$ws.Range("E28").Value = "  +6.07%  "
$ws.Range("D29").Value = "'0.113"
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "'3.86"
$ws.Range("E31").Value = "  +10.52%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.0512"
$ws.Range("E32").Value = "  +2.98%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.20"
$ws.Range("E33").Value = "  +5.18%  "
$ws.Range("B34").Value = "MinaProtocolToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
$ws.Range("D34").Value = "'1.69"
$ws.Range("E34").Value = "  +309.57%  "
$ws.Range("D35").Value = "'3.51"
$ws.Range("E35").Value = "  +6.14%  "
$ws.Range("D36").Value = "'1.538.79"
$ws.Range("E36").Value = "  +1.56%  "
$ws.Range("E37").Value = "  +2.57%  "
$ws.Range("D38").Value = "'1.08"
$ws.Range("E38").Value = "  +4.59%  "
$ws.Range("D39").Value = "'84.61"
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("E40").Value = "  +4.55%  "
$ws.Range("D41").Value = "'0.623"
$ws.Range("E41").Value = "  +4.56%  "
$ws.Range("E42").Value = "  +2.88%  "
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("D44").Value = "'0.910"
$ws.Range("E44").Value = "  +8.13%  "
$ws.Range("E45").Value = "  +7.74%  "
$ws.Range("E46").Value = "  +4.10%  "
$ws.Range("E47").Value = "  +4.32%  "
$ws.Range("D48").Value = "'1.970.54"
$ws.Range("E48").Value = "  +8.39%  "
$ws.Range("E49").Value = "  +2.75%  "
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").Value = "'52.76"
$ws.Range("E51").Value = "  +1.68%  "
